$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Done" status column (I) - the app no longer tracks a Done flag per image size.
$ws.Range("I2:I7").ClearContents()

# Logo width (F) now matches the background width (E) instead of being slightly narrower -
# the logo is no longer centered with an offset, so G (Logo top) recalculates automatically.
$ws.Range("F2:F7").Value = 620

# Drop the two custom resolution rows (1024x600 / 600x1024) - only the standard presets remain.
$ws.Rows("8:9").Delete()

# Re-establish the shared formulas over the now-shorter range (rows 3-7).
$ws.Range("C3:C7").Formula = "=B3/3"
$ws.Range("D3:D7").Formula = "=2*B3/3"

# Update the active selection to the resized Logo V column.
[void]$ws.Range("F2:F7").Select()
